$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 currently holds the text "30" (inline string); change it to the numeric value 30
$ws.Range("B4").Value = 30

# Add a new row 5 with data: A5="test", B5="30" (kept as text, not a number), C5="hi"
$ws.Range("A5").Value = "test"

# Force B5 to be stored as text "30" (not auto-converted to a number),
# then restore its style to the default so no extra formatting sticks.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "30"
$ws.Range("B5").Style = $ws.Range("A5").Style

$ws.Range("C5").Value = "hi"
